$wb = $excel.ActiveWorkbook

# --- Rename existing sheet to "Prompts" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Prompts"

# --- Add second sheet "Semantics" after Prompts ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Semantics"

# =========================================================
# Sheet "Prompts" data
# =========================================================
$ws1.Range("A1").Value = "TestCaseName"
$ws1.Range("B1").Value = "ConnectionName"
$ws1.Range("C1").Value = "Prompt1"
$ws1.Range("D1").Value = "Prompt2"

$ws1.Range("A2").Value = "test_verifyInsurancePrompt1"
$ws1.Range("B2").Value = "cementdemo"
$ws1.Range("C2").Value = "What are the sales record for past year"
$ws1.Range("D2").Value = "What are the sales record for past quarter"

$ws1.Range("A3").Value = "test_verifyInsurancePrompt2"
$ws1.Range("B3").Value = "demo_retail"
$ws1.Range("C3").Value = "Show me the Total Net Sales by Month with units sold and Year."

$ws1.Range("A4").Value = "test_verifyInsurancePrompt3"
$ws1.Range("B4").Value = "demo_retail"
$ws1.Range("C4").Value = "Show me Net Sales and Units Sold by Category and subcategory."

$ws1.Range("A5").Value = "test_verifyInsurancePrompt4"
$ws1.Range("B5").Value = "demo_retail"
$ws1.Range("C5").Value = "Show me Net Sales and Units Sold by Category and subcategory against month and year."
$ws1.Range("D5").Value = "filter for the month of march."

$ws1.Range("A6").Value = "test_verifyInsurancePrompt5"
$ws1.Range("B6").Value = "demo_retail"
$ws1.Range("C6").Value = "Identify our top 10 customers by calculating the Average Order Value (AOV = AVG(SalesAmount))"

$ws1.Range("A7").Value = "test_verifyInsurancePrompt6"
$ws1.Range("B7").Value = "demo_retail"
$ws1.Range("C7").Value = "Show me the Year-to-Date total Net Sales and Quantity Sold by Store"

$ws1.Range("A8").Value = "test_verifyInsurancePrompt7"
$ws1.Range("B8").Value = "demo_retail"
$ws1.Range("C8").Value = "Compare sales of each product based on discount (DiscountPercent vs Net Sales)"

$ws1.Range("A9").Value = "test_verifyInsurancePrompt8"
$ws1.Range("B9").Value = "demo_retail"
$ws1.Range("C9").Value = "Show Profit and Margin % by Category"

$ws1.Range("A10").Value = "test_verifyInsurancePrompt9"
$ws1.Range("B10").Value = "demo_retail"
$ws1.Range("C10").Value = "Show Profit and Margin % by Category and subcategory"

# column widths on Prompts sheet (nearest values this engine's pixel-quantized
# column-width storage can represent to the authored widths 78.33203125 / 34.5546875)
$ws1.Columns.Item(3).ColumnWidth = 77.5
$ws1.Columns.Item(4).ColumnWidth = 33.666666666666664

# selection on Prompts sheet
$ws1.Range("D2").Select() | Out-Null

# =========================================================
# Sheet "Semantics" data
# =========================================================
$ws2.Range("A1").Value = "TestCaseName"
$ws2.Range("B1").Value = "ConnectionName"
$ws2.Range("C1").Value = "Semantic Views"

$ws2.Range("A2").Value = "test_verifyInsuranceSemantic1"
$ws2.Range("B2").Value = "demo_retail"
$ws2.Range("C2").Value = "SalesByQty"

$ws2.Range("A3").Value = "test_verifyInsuranceSemantic2"
$ws2.Range("B3").Value = "demo_retail"
$ws2.Range("C3").Value = "discount vs sales"

$ws2.Range("A4").Value = "test_verifyInsuranceSemantic3"
$ws2.Range("B4").Value = "demo_retail"
$ws2.Range("C4").Value = "YTD Net Sales"

$ws2.Range("A5").Value = "test_verifyInsuranceSemantic4"
$ws2.Range("B5").Value = "demo_retail"
$ws2.Range("C5").Value = "Sales Comparision"

# column widths on Semantics sheet (nearest representable values to the
# authored widths 26 / 15.21875 / 14.5546875)
$ws2.Columns.Item(1).ColumnWidth = 25.166666666666668
$ws2.Columns.Item(2).ColumnWidth = 14.333333333333332
$ws2.Columns.Item(3).ColumnWidth = 13.666666666666666

# Semantics is the active/visible tab
$ws2.Activate()
